$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1523.8
$ws.Range("I20").Value = 1155
$ws.Range("K20").Value = 1155
$ws.Range("M20").Value = -925

$ws.Range("H31").Value = 340.77777
$ws.Range("I31").Value = 352.125
$ws.Range("K31").Value = 1056.375
$ws.Range("M31").Value = -826.375

$ws.Range("H35").Value = 1523.8
$ws.Range("I35").Value = 1155
$ws.Range("K35").Value = 1155
$ws.Range("M35").Value = -776

$ws.Range("H41").Value = 858.26666
$ws.Range("J41").Value = 179.8
$ws.Range("L41").Value = 179.8
$ws.Range("N41").Value = -1059.8

$ws.Range("H51").Value = 16670265
$ws.Range("I51").Value = 3995
$ws.Range("J51").Value = 27781112
$ws.Range("K51").Value = 3995
$ws.Range("L51").Value = 27781112
$ws.Range("M51").Value = -3511
$ws.Range("N51").Value = -27782080

$ws.Range("H106").Value = 17853.572
$ws.Range("I106").Value = 3246.25
$ws.Range("J106").Value = 37330
$ws.Range("K106").Value = 3246.25
$ws.Range("L106").Value = 37330
$ws.Range("M106").Value = -2615.25
$ws.Range("N106").Value = -38592

$ws.Range("H131").Value = 735.8889
$ws.Range("I131").Value = 678
$ws.Range("K131").Value = 2034
$ws.Range("M131").Value = 3006

$ws.Range("H138").Value = 2868.36
$ws.Range("I138").Value = 4912.6665
$ws.Range("J138").Value = 2222.7896
$ws.Range("K138").Value = 14737.9995
$ws.Range("L138").Value = 6668.3688
$ws.Range("M138").Value = -9597.999500000002
$ws.Range("N138").Value = -16948.3688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 62.1
$ws.Range("I5").Value = 52.875
$ws.Range("J5").Value = 99
$ws.Range("K5").Value = 52.875
$ws.Range("L5").Value = 99
$ws.Range("M5").Value = 59.125
$ws.Range("N5").Value = -323

$ws.Range("H32").Value = 30450.906
$ws.Range("I32").Value = 6969.3687
$ws.Range("J32").Value = 208910.6
$ws.Range("K32").Value = 6969.3687
$ws.Range("L32").Value = 208910.6
$ws.Range("M32").Value = -6682.3687
$ws.Range("N32").Value = -209484.6

$ws.Range("H102").Value = 3192
$ws.Range("I102").Value = 1366.7142
$ws.Range("K102").Value = 1366.7142
$ws.Range("M102").Value = 255.2858000000001

$ws.Range("H110").Value = 2227.9285
$ws.Range("I110").Value = 2290.0908
$ws.Range("K110").Value = 2290.0908
$ws.Range("M110").Value = -245.0907999999999

$ws.Range("H112").Value = 35521.832
$ws.Range("J112").Value = 35521.832
$ws.Range("L112").Value = 35521.832
$ws.Range("N112").Value = -38475.832

$ws.Range("H122").Value = 1930
$ws.Range("I122").Value = 1847.5358
$ws.Range("J122").Value = 2314.8333
$ws.Range("K122").Value = 5542.607400000001
$ws.Range("L122").Value = 6944.499899999999
$ws.Range("M122").Value = -3092.607400000001
$ws.Range("N122").Value = -11844.4999

$ws.Range("H132").Value = 2348.4827
$ws.Range("I132").Value = 913.9091
$ws.Range("K132").Value = 2741.7273
$ws.Range("M132").Value = -211.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 62.1
$ws.Range("I4").Value = 52.875
$ws.Range("J4").Value = 99
$ws.Range("K4").Value = 52.875
$ws.Range("L4").Value = 99
$ws.Range("M4").Value = 62.125
$ws.Range("N4").Value = -329

$ws.Range("H105").Value = 25405
$ws.Range("I105").Value = 50010
$ws.Range("J105").Value = 800
$ws.Range("K105").Value = 50010
$ws.Range("L105").Value = 800
$ws.Range("M105").Value = -48263
$ws.Range("N105").Value = -4294

$ws.Range("H107").Value = 67911.47
$ws.Range("I107").Value = 112313.78
$ws.Range("J107").Value = 1308
$ws.Range("K107").Value = 112313.78
$ws.Range("L107").Value = 1308
$ws.Range("M107").Value = -110393.78
$ws.Range("N107").Value = -5148

$ws.Range("H108").Value = 74421
$ws.Range("J108").Value = 74421
$ws.Range("L108").Value = 74421
$ws.Range("N108").Value = -82101

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1841.6111
$ws.Range("I58").Value = 1819.0869
$ws.Range("K58").Value = 1819.0869
$ws.Range("M58").Value = -1616.0869

$ws.Range("H94").Value = 2067
$ws.Range("J94").Value = 1065
$ws.Range("L94").Value = 1065
$ws.Range("N94").Value = -1967

$ws.Range("H99").Value = 31025.357
$ws.Range("I99").Value = 31565.77
$ws.Range("K99").Value = 31565.77
$ws.Range("M99").Value = -30067.77

$ws.Range("H107").Value = 1538.64
$ws.Range("I107").Value = 1336.8
$ws.Range("J107").Value = 1841.4
$ws.Range("K107").Value = 1336.8
$ws.Range("L107").Value = 1841.4
$ws.Range("M107").Value = 583.2
$ws.Range("N107").Value = -5681.4

$ws.Range("H126").Value = 31025.357
$ws.Range("I126").Value = 31565.77
$ws.Range("K126").Value = 94697.31
$ws.Range("M126").Value = -92227.31

$ws.Range("H136").Value = 1841.6111
$ws.Range("I136").Value = 1819.0869
$ws.Range("K136").Value = 5457.2607
$ws.Range("M136").Value = -2907.2607

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 406.75
$ws.Range("J5").Value = 297
$ws.Range("L5").Value = 891
$ws.Range("N5").Value = -1115

$ws.Range("H29").Value = 53.833332
$ws.Range("I29").Value = 47
$ws.Range("J29").Value = 67.5
$ws.Range("K29").Value = 141
$ws.Range("L29").Value = 202.5
$ws.Range("M29").Value = 136
$ws.Range("N29").Value = -756.5

$ws.Range("H70").Value = 3248
$ws.Range("I70").Value = 3248
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 9744
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -9429
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 3248
$ws.Range("I73").Value = 3248
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 9744
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8652
$ws.Range("N73").ClearContents()

$ws.Range("H131").Value = 2921.6667
$ws.Range("J131").Value = 2875
$ws.Range("L131").Value = 8625
$ws.Range("N131").Value = -18705

$ws.Range("H135").Value = 406.75
$ws.Range("J135").Value = 297
$ws.Range("L135").Value = 2673
$ws.Range("N135").Value = -7743

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 22399.871
$ws.Range("I97").Value = 30999.727
$ws.Range("J97").Value = 1378
$ws.Range("K97").Value = 30999.727
$ws.Range("L97").Value = 1378
$ws.Range("M97").Value = -30503.727
$ws.Range("N97").Value = -2370

$ws.Range("H132").Value = 3405.4075
$ws.Range("I132").Value = 2731
$ws.Range("K132").Value = 8193
$ws.Range("M132").Value = -5663

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 2999.3333
$ws.Range("J34").Value = 2999.5
$ws.Range("L34").Value = 2999.5
$ws.Range("N34").Value = -3343.5

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H110").Value = 57881.332
$ws.Range("J110").Value = 57881.332
$ws.Range("L110").Value = 57881.332
$ws.Range("N110").Value = -66061.33199999999

$ws.Range("H136").Value = 4352.92
$ws.Range("I136").Value = 3965.5
$ws.Range("K136").Value = 11896.5
$ws.Range("M136").Value = -9346.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1104.1
$ws.Range("I13").Value = 1551
$ws.Range("J13").Value = 657.2
$ws.Range("K13").Value = 1551
$ws.Range("L13").Value = 657.2
$ws.Range("M13").Value = -1411
$ws.Range("N13").Value = -937.2

$ws.Range("H70").Value = 20105
$ws.Range("J70").Value = 20105
$ws.Range("L70").Value = 20105
$ws.Range("N70").Value = -20735

$ws.Range("H73").Value = 20105
$ws.Range("J73").Value = 20105
$ws.Range("L73").Value = 20105
$ws.Range("N73").Value = -22289

$ws.Range("H122").Value = 884.5217
$ws.Range("I122").Value = 901.9048
$ws.Range("J122").Value = 702
$ws.Range("K122").Value = 2705.7144
$ws.Range("L122").Value = 2106
$ws.Range("M122").Value = -255.7143999999998
$ws.Range("N122").Value = -7006
